$d = $word.ActiveDocument

# 1. Fix the post-command example: ("Group2") -> ("Group0")
#    (search/replace only the digit so the surrounding straight quotes are
#    left completely untouched and are not smart-quoted by AutoFormat)
$found1 = $d.Content.Find.Execute(
    "Group2",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Group0",
    2)
Write-Host "Group2 -> Group0 replaced:" $found1
if (-not $found1) {
    throw "Could not find 'Group2' to replace with 'Group0'"
}

# 2. Remove the three superseded `pt_base` example lines (age / qol / gender
#    type(bin)) together with the manual line breaks that introduced them,
#    leaving the single break that already separated the surrounding lines.
$search2 = "^l      . pt_base age , post(``postname') over(treat)  overall(last) over_grps(1, 0) type(cont) su_label(append)^l      . pt_base qol , post(``postname') over(treat)  overall(last)  over_grps(1, 0) type(skew) su_label(append)^l      . pt_base  gender , post(``postname') over(treat)  overall(last)  over_grps(1, 0) type(bin)  su_label(append)"
$found2 = $d.Content.Find.Execute(
    $search2,
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "",
    2)
Write-Host "pt_base age/qol/gender(bin) lines removed:" $found2
if (-not $found2) {
    throw "Could not find the pt_base age/qol/gender(bin) example block to remove"
}
